# "combined marraige and regular data"
#
# The workbook's data rows (2:10) were selected as a block (e.g. after
# copying/combining data from another sheet) and the data row height was
# tightened to 10pt.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the full data rows (2 through 10, whole rows to the end of the grid)
# -> sheetView selection activeCell="A2" sqref="A2:XFD10"
$ws.Range("A2:XFD10").Select()

# Row 2 gets a custom (smaller) height
# -> row r="2" ht="10" customHeight="1"
$ws.Rows("2:2").RowHeight = 10
